# Auto-generated Excel COM-interop script
# Applies scheduled-runner market price/profit updates to Anima_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Range("H43").Value = 2033.3334
$ws.Range("J43").Value = 2033.3334
$ws.Range("L43").Value = 2033.3334
$ws.Range("N43").Value = -2171.3334

# Row 113: Amaro Kart | Starch Glue
$ws.Range("H113").Value = 3354.2856
$ws.Range("I113").Value = 2966.6667
$ws.Range("J113").Value = 3645
$ws.Range("K113").Value = 2966.6667
$ws.Range("L113").Value = 3645
$ws.Range("M113").Value = 287.3332999999998
$ws.Range("N113").Value = -10153

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Range("H137").Value = 2951.4614
$ws.Range("I137").Value = 2422
$ws.Range("J137").Value = 5371.857
$ws.Range("K137").Value = 7266
$ws.Range("L137").Value = 16115.571
$ws.Range("M137").Value = -4716
$ws.Range("N137").Value = -21215.571

$ws = $wb.Worksheets.Item("ARM")
# Row 41: Skillet Scandal | White Skillet
$ws.Range("H41").Value = 1514.25
$ws.Range("I41").Value = 1514.25
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1514.25
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1100.25
$ws.Range("N41").ClearContents()

# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2398.279
$ws.Range("I61").Value = 1872.6666
$ws.Range("J61").Value = 4132.8
$ws.Range("K61").Value = 1872.6666
$ws.Range("L61").Value = 4132.8
$ws.Range("M61").Value = -1660.6666
$ws.Range("N61").Value = -4556.8

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1227.1538
$ws.Range("I74").Value = 730.2143
$ws.Range("K74").Value = 730.2143
$ws.Range("M74").Value = 143.7857

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1227.1538
$ws.Range("I77").Value = 730.2143
$ws.Range("K77").Value = 3651.0715
$ws.Range("M77").Value = 716.9285

# Row 97: Ore for Me | High Steel Ingot
$ws.Range("H97").Value = 921.95654
$ws.Range("I97").Value = 831.25
$ws.Range("J97").Value = 1129.2858
$ws.Range("K97").Value = 831.25
$ws.Range("L97").Value = 1129.2858
$ws.Range("M97").Value = -335.25
$ws.Range("N97").Value = -2121.2858

# Row 110: Scheduled Maintenance | Deepgold Ingot
$ws.Range("H110").Value = 1129
$ws.Range("I110").Value = 1136.5
$ws.Range("J110").Value = 1116.1428
$ws.Range("K110").Value = 1136.5
$ws.Range("L110").Value = 1116.1428
$ws.Range("M110").Value = 908.5
$ws.Range("N110").Value = -5206.1428

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2398.279
$ws.Range("I136").Value = 1872.6666
$ws.Range("J136").Value = 4132.8
$ws.Range("K136").Value = 5617.9998
$ws.Range("L136").Value = 12398.4
$ws.Range("M136").Value = -3067.9998
$ws.Range("N136").Value = -17498.4

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 1398.0714
$ws.Range("I94").Value = 1058.8572
$ws.Range("J94").Value = 2415.7144
$ws.Range("K94").Value = 1058.8572
$ws.Range("L94").Value = 2415.7144
$ws.Range("M94").Value = -607.8571999999999
$ws.Range("N94").Value = -3317.7144

# Row 132: Always Be Prepaired | Mountain Chromite Twinfangs
$ws.Range("H132").Value = 42158.234
$ws.Range("J132").Value = 42158.234
$ws.Range("L132").Value = 42158.234
$ws.Range("N132").Value = -52278.234

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 1955.775
$ws.Range("I134").Value = 1522.1515
$ws.Range("K134").Value = 4566.4545
$ws.Range("M134").Value = -2031.4545

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 333.33334
$ws.Range("J7").Value = 200
$ws.Range("L7").Value = 200
$ws.Range("N7").Value = -426

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 5966.2
$ws.Range("I31").Value = 1394.6296
$ws.Range("J31").Value = 11332.826
$ws.Range("K31").Value = 1394.6296
$ws.Range("L31").Value = 11332.826
$ws.Range("M31").Value = -1099.6296
$ws.Range("N31").Value = -11922.826

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 5966.2
$ws.Range("I34").Value = 1394.6296
$ws.Range("J34").Value = 11332.826
$ws.Range("K34").Value = 1394.6296
$ws.Range("L34").Value = 11332.826
$ws.Range("M34").Value = -1192.6296
$ws.Range("N34").Value = -11736.826

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 3624845
$ws.Range("I132").Value = 1166.3334
$ws.Range("K132").Value = 3499.0002
$ws.Range("M132").Value = -969.0001999999999

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 4566.5835
$ws.Range("I134").Value = 5092.96
$ws.Range("J134").Value = 3370.2727
$ws.Range("K134").Value = 15278.88
$ws.Range("L134").Value = 10110.8181
$ws.Range("M134").Value = -12743.88
$ws.Range("N134").Value = -15180.8181

$ws = $wb.Worksheets.Item("CUL")
# Row 9: Jack of All Plates | Jack-o'-lantern
$ws.Range("H9").Value = 112987.625
$ws.Range("J9").Value = 129057.29
$ws.Range("L9").Value = 387171.87
$ws.Range("N9").Value = -387619.87

# Row 17: Chew the Fat | Grilled Dodo
$ws.Range("H17").Value = 634.625
$ws.Range("J17").Value = 969.25
$ws.Range("L17").Value = 2907.75
$ws.Range("N17").Value = -3245.75

# Row 20: Omelette's Be Friends | Dodo Omelette
$ws.Range("H20").Value = 1359.9
$ws.Range("I20").Value = 999
$ws.Range("J20").Value = 1400
$ws.Range("K20").Value = 2997
$ws.Range("L20").Value = 4200
$ws.Range("M20").Value = -2770
$ws.Range("N20").Value = -4654

# Row 22: A Total Nut Job | Walnut Bread
$ws.Range("H22").Value = 999.8889
$ws.Range("I22").Value = 999.8
$ws.Range("K22").Value = 2999.4
$ws.Range("M22").Value = -2830.4

# Row 26: A Grape Idea | Grape Juice
$ws.Range("H26").Value = 376.4737
$ws.Range("I26").Value = 75.666664
$ws.Range("J26").Value = 515.3077
$ws.Range("K26").Value = 226.999992
$ws.Range("L26").Value = 1545.9231
$ws.Range("M26").Value = 61.00000800000001
$ws.Range("N26").Value = -2121.9231

# Row 27: Brain Food | Walnut Bread
$ws.Range("H27").Value = 999.8889
$ws.Range("I27").Value = 999.8
$ws.Range("K27").Value = 2999.4
$ws.Range("M27").Value = -2897.4

# Row 32: Convalescence Precedes Essence | Ginger Cookie
$ws.Range("H32").Value = 11114258
$ws.Range("J32").Value = 11114258
$ws.Range("L32").Value = 33342774
$ws.Range("N32").Value = -33343340

# Row 34: Fever Pitch | Chamomile Tea
$ws.Range("H34").Value = 11364130
$ws.Range("J34").Value = 11905270
$ws.Range("L34").Value = 35715810
$ws.Range("N34").Value = -35715978

# Row 39: Bloody Good Tart, This | Blood Currant Tart
$ws.Range("H39").Value = 2980
$ws.Range("J39").Value = 2980
$ws.Range("L39").Value = 8940
$ws.Range("N39").Value = -9528

# Row 51: The Perks of Life at Sea | Jerked Beef
$ws.Range("H51").Value = 1212.6666
$ws.Range("I51").Value = 800
$ws.Range("J51").Value = 1242.1428
$ws.Range("K51").Value = 2400
$ws.Range("L51").Value = 3726.4284
$ws.Range("M51").Value = -1940
$ws.Range("N51").Value = -4646.428400000001

# Row 57: The Egg Files | Deviled Eggs
$ws.Range("H57").Value = 1493.9584
$ws.Range("I57").Value = 501
$ws.Range("J57").Value = 1755.2632
$ws.Range("K57").Value = 1503
$ws.Range("L57").Value = 5265.7896
$ws.Range("M57").Value = -944
$ws.Range("N57").Value = -6383.7896

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1031.0435
$ws.Range("I131").Value = 631.2857
$ws.Range("J131").Value = 1205.9375
$ws.Range("K131").Value = 1893.8571
$ws.Range("L131").Value = 3617.8125
$ws.Range("M131").Value = 3146.1429
$ws.Range("N131").Value = -13697.8125

$ws = $wb.Worksheets.Item("GSM")
# Row 48: Dead Can't Defang | Wolf Necklace
$ws.Range("H48").Value = 6500
$ws.Range("I48").Value = 3000
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 3000
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = -2515
$ws.Range("N48").Value = -10970

# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 2544.4075
$ws.Range("I132").Value = 2160.6086
$ws.Range("J132").Value = 4751.25
$ws.Range("K132").Value = 6481.825800000001
$ws.Range("L132").Value = 14253.75
$ws.Range("M132").Value = -3951.825800000001
$ws.Range("N132").Value = -19313.75

$ws = $wb.Worksheets.Item("LTW")
# Row 64: Glorified Hole-punchers | Archaeoskin Gloves of Aiming
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67: Treat Them with Kid Gloves (L) | Archaeoskin Gloves of Aiming
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 64: Ribbon of Remembrance | Rainbow Ribbon of Healing
$ws.Range("H64").Value = 27500
$ws.Range("J64").Value = 27500
$ws.Range("L64").Value = 27500
$ws.Range("N64").Value = -27996

# Row 67: The Road Was a Ribbon of Moonlight (L) | Rainbow Ribbon of Healing
$ws.Range("H67").Value = 27500
$ws.Range("J67").Value = 27500
$ws.Range("L67").Value = 27500
$ws.Range("N67").Value = -29216

# Row 74: Clothing the Naked Truth | Ramie Robe of Casting
$ws.Range("H74").Value = 10400
$ws.Range("J74").Value = 10600
$ws.Range("L74").Value = 10600
$ws.Range("N74").Value = -12472

# Row 77: When in Robes (L) | Ramie Robe of Casting
$ws.Range("H77").Value = 10400
$ws.Range("J77").Value = 10600
$ws.Range("L77").Value = 31800
$ws.Range("N77").Value = -41160

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2302.72
$ws.Range("I122").Value = 2258.7778
$ws.Range("J122").Value = 2415.7144
$ws.Range("K122").Value = 6776.3334
$ws.Range("L122").Value = 7247.1432
$ws.Range("M122").Value = -4326.3334
$ws.Range("N122").Value = -12147.1432

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 5749816
$ws.Range("I132").Value = 2978.2942
$ws.Range("J132").Value = 13891170
$ws.Range("K132").Value = 8934.882599999999
$ws.Range("L132").Value = 41673510
$ws.Range("M132").Value = -6404.882599999999
$ws.Range("N132").Value = -41678570

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 2627.1892
$ws.Range("I136").Value = 2115.037
$ws.Range("K136").Value = 6345.110999999999
$ws.Range("M136").Value = -3795.110999999999
